$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (header is in row 1). The commit bumps that date by one day
# (45180 -> 45181) for every populated row.
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current + 1
    }
}
